# Adding MPA test automation upload file
# Update the auto-generated transfer/upload reference numbers on the "Data" sheet.
# The sequential IDs shift forward (e.g. 60000373 -> 60000391, 271 -> 278,
# 60000374 -> 60000392, 272 -> 279) for every data row (6-29).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("K6").Value = 60000391

$ws.Range("L7").Value = 278
$ws.Range("N7").Value = 60000392

$ws.Range("L8").Value = 278
$ws.Range("O8").Value = 279

$ws.Range("K9").Value = 60000391

$ws.Range("L10").Value = 278

$ws.Range("K11").Value = 60000391
$ws.Range("N11").Value = 60000392

$ws.Range("L12").Value = 278
$ws.Range("N12").Value = 60000392

$ws.Range("L13").Value = 278
$ws.Range("O13").Value = 279

$ws.Range("K14").Value = 60000391

$ws.Range("L15").Value = 278

$ws.Range("K16").Value = 60000391
$ws.Range("N16").Value = 60000392

$ws.Range("L17").Value = 278
$ws.Range("N17").Value = 60000392

$ws.Range("L18").Value = 278
$ws.Range("O18").Value = 279

$ws.Range("K19").Value = 60000391

$ws.Range("L20").Value = 278

$ws.Range("K21").Value = 60000391
$ws.Range("N21").Value = 60000392

$ws.Range("L22").Value = 278
$ws.Range("N22").Value = 60000392

$ws.Range("L23").Value = 278
$ws.Range("O23").Value = 279

$ws.Range("K24").Value = 60000391

$ws.Range("L25").Value = 278

$ws.Range("K26").Value = 60000391
$ws.Range("N26").Value = 60000392

$ws.Range("L27").Value = 278
$ws.Range("N27").Value = 60000392

$ws.Range("L28").Value = 278
$ws.Range("O28").Value = 279

$ws.Range("K29").Value = 60000391

$wb.Save()
